# Auto update Excel log
# Appends new sensor-log rows to the "PIR" sheet (rows 238-251) and the
# "Humidity" sheet (rows 157-169), matching the source data export.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Writes the data rows into the given worksheet starting at startRow.
# Column A ("Date") would otherwise be auto-interpreted by Excel as a
# real date, and (when requested) column E ("Value") would be
# auto-interpreted as a numeric percentage. Both are briefly switched to
# Text format so the literal string is preserved exactly, then the
# formatting is cleared again so the new cells keep the workbook's
# default ("General") style, just like the existing rows.
# ---------------------------------------------------------------------
function Write-LogRows($ws, $startRow, $rows, $valueColumnIsText) {
    $rowCount = $rows.Count
    $lastRow = $startRow + $rowCount - 1

    $dateRange = $ws.Range("A" + $startRow + ":A" + $lastRow)
    $dateRange.NumberFormat = "@"

    if ($valueColumnIsText) {
        $valueRange = $ws.Range("E" + $startRow + ":E" + $lastRow)
        $valueRange.NumberFormat = "@"
    }

    for ($i = 0; $i -lt $rowCount; $i++) {
        $r = $startRow + $i
        $data = $rows[$i]
        $ws.Cells.Item($r, 1).Value = $data[0]
        $ws.Cells.Item($r, 2).Value = $data[1]
        $ws.Cells.Item($r, 3).Value = $data[2]
        $ws.Cells.Item($r, 4).Value = $data[3]
        $ws.Cells.Item($r, 5).Value = $data[4]
        $ws.Cells.Item($r, 6).Value = $data[5]
    }

    $dateRange.ClearFormats()
    if ($valueColumnIsText) {
        $valueRange.ClearFormats()
    }
}

# ---------------------------------------------------------------------
# PIR sheet new rows (Date, Timestamp, Hour, Location, Value, Status)
# ---------------------------------------------------------------------
$pirRows = @(
    @("2026-01-30","17:26:48","17:00","Bathroom","No Motion","Inactive"),
    @("2026-01-30","17:26:49","17:00","Bathroom","No Motion","Inactive"),
    @("2026-01-30","17:26:53","17:00","Bathroom","No Motion","Inactive"),
    @("2026-01-30","17:26:58","17:00","Bathroom","No Motion","Inactive"),
    @("2026-01-30","17:27:03","17:00","Bathroom","No Motion","Inactive"),
    @("2026-01-30","17:27:08","17:00","Bathroom","No Motion","Inactive"),
    @("2026-01-30","17:27:13","17:00","Bathroom","No Motion","Inactive"),
    @("2026-01-30","17:27:18","17:00","Bathroom","No Motion","Inactive"),
    @("2026-01-30","17:27:23","17:00","Bathroom","No Motion","Inactive"),
    @("2026-01-30","17:27:28","17:00","Bathroom","No Motion","Inactive"),
    @("2026-01-30","17:27:33","17:00","Bathroom","No Motion","Inactive"),
    @("2026-01-30","17:27:38","17:00","Bathroom","No Motion","Inactive"),
    @("2026-01-30","17:27:43","17:00","Bathroom","No Motion","Inactive"),
    @("2026-01-30","17:27:48","17:00","Bathroom","No Motion","Inactive")
)

$wsPir = $wb.Worksheets.Item("PIR")
Write-LogRows $wsPir 238 $pirRows $false

# ---------------------------------------------------------------------
# Humidity sheet new rows (Date, Timestamp, Hour, Location, Value, Status)
# ---------------------------------------------------------------------
$humidityRows = @(
    @("2026-01-30","17:26:49","17:00","Bathroom","87.3%","Active"),
    @("2026-01-30","17:26:50","17:00","Bathroom","87.4%","Active"),
    @("2026-01-30","17:26:54","17:00","Bathroom","87.3%","Active"),
    @("2026-01-30","17:26:59","17:00","Bathroom","87.4%","Active"),
    @("2026-01-30","17:27:04","17:00","Bathroom","86.4%","Active"),
    @("2026-01-30","17:27:09","17:00","Bathroom","87.4%","Active"),
    @("2026-01-30","17:27:14","17:00","Bathroom","87.4%","Active"),
    @("2026-01-30","17:27:19","17:00","Bathroom","87.4%","Active"),
    @("2026-01-30","17:27:24","17:00","Bathroom","87.4%","Active"),
    @("2026-01-30","17:27:29","17:00","Bathroom","87.4%","Active"),
    @("2026-01-30","17:27:34","17:00","Bathroom","87.4%","Active"),
    @("2026-01-30","17:27:39","17:00","Bathroom","87.4%","Active"),
    @("2026-01-30","17:27:44","17:00","Bathroom","86.4%","Active")
)

$wsHumidity = $wb.Worksheets.Item("Humidity")
Write-LogRows $wsHumidity 157 $humidityRows $true
